$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

$ws.Cells.Item($row, 1).Value = "r585"
$ws.Cells.Item($row, 2).Value = "infinite loop"
$ws.Cells.Item($row, 3).Value = "fixed"
$ws.Cells.Item($row, 4).Value = "2025-10-01 14:50:45"
